# Auto-generated edit script: applies numeric updates to Hades_Profits-style sheets
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 33
$ws.Range("H33").Value = 256.5
$ws.Range("I33").Value = 134.2
$ws.Range("J33").Value = 460.33334
$ws.Range("K33").Value = 134.2
$ws.Range("L33").Value = 460.33334
$ws.Range("M33").Value = 94.80000000000001
$ws.Range("N33").Value = -918.33334
# Row 97
$ws.Range("H97").Value = 1040.8334
$ws.Range("J97").Value = 1040.8334
$ws.Range("L97").Value = 3122.5002
$ws.Range("N97").Value = -4114.5002
# Row 100
$ws.Range("H100").Value = 1316.6666
$ws.Range("I100").Value = 1360
$ws.Range("K100").Value = 1360
$ws.Range("M100").Value = -819
# Row 103
$ws.Range("H103").Value = 657.5
$ws.Range("I103").Value = 716.6667
$ws.Range("J103").Value = 598.3333
$ws.Range("K103").Value = 2150.0001
$ws.Range("L103").Value = 1794.9999
$ws.Range("M103").Value = -1564.0001
$ws.Range("N103").Value = -2966.9999
# Row 106
$ws.Range("H106").Value = 4800.5557
$ws.Range("I106").Value = 4800.5557
$ws.Range("J106").Value = 0
$ws.Range("K106").Value = 4800.5557
$ws.Range("L106").Value = 0
$ws.Range("M106").Value = -4169.5557
$ws.Range("N106").ClearContents()
# Row 113
$ws.Range("H113").Value = 3647
$ws.Range("I113").Value = 4015
$ws.Range("J113").Value = 3325
$ws.Range("K113").Value = 4015
$ws.Range("L113").Value = 3325
$ws.Range("M113").Value = -761
$ws.Range("N113").Value = -9833
# Row 127
$ws.Range("H127").Value = 1092.05
$ws.Range("I127").Value = 484.55554
$ws.Range("K127").Value = 1453.66662
$ws.Range("M127").Value = 3506.33338

$ws = $wb.Worksheets.Item("ARM")
# Row 2
$ws.Range("H2").Value = 1302.7778
$ws.Range("I2").Value = 1276.7391
$ws.Range("J2").Value = 1452.5
$ws.Range("K2").Value = 1276.7391
$ws.Range("L2").Value = 1452.5
$ws.Range("M2").Value = -1163.7391
$ws.Range("N2").Value = -1678.5
# Row 45
$ws.Range("H45").Value = 1683.8235
$ws.Range("I45").Value = 1790.5518
$ws.Range("J45").Value = 1064.8
$ws.Range("K45").Value = 1790.5518
$ws.Range("L45").Value = 1064.8
$ws.Range("M45").Value = -1413.5518
$ws.Range("N45").Value = -1818.8
# Row 102
$ws.Range("H102").Value = 1607.6428
$ws.Range("I102").Value = 1411.8
$ws.Range("K102").Value = 1411.8
$ws.Range("M102").Value = 210.2
# Row 116
$ws.Range("H116").Value = 1302.7778
$ws.Range("I116").Value = 1276.7391
$ws.Range("J116").Value = 1452.5
$ws.Range("K116").Value = 1276.7391
$ws.Range("L116").Value = 1452.5
$ws.Range("M116").Value = 1017.2609
$ws.Range("N116").Value = -6040.5
# Row 132
$ws.Range("H132").Value = 53633.49
$ws.Range("I132").Value = 40549.73
$ws.Range("J132").Value = 79801
$ws.Range("K132").Value = 121649.19
$ws.Range("L132").Value = 239403
$ws.Range("M132").Value = -119119.19
$ws.Range("N132").Value = -244463

$ws = $wb.Worksheets.Item("BSM")
# Row 3
$ws.Range("H3").Value = 1302.7778
$ws.Range("I3").Value = 1276.7391
$ws.Range("J3").Value = 1452.5
$ws.Range("K3").Value = 1276.7391
$ws.Range("L3").Value = 1452.5
$ws.Range("M3").Value = -1162.7391
$ws.Range("N3").Value = -1680.5
# Row 99
$ws.Range("H99").Value = 1316.1333
$ws.Range("I99").Value = 1251.125
$ws.Range("J99").Value = 1390.4286
$ws.Range("K99").Value = 1251.125
$ws.Range("L99").Value = 1390.4286
$ws.Range("M99").Value = 246.875
$ws.Range("N99").Value = -4386.4286
# Row 134
$ws.Range("H134").Value = 1850.9584
$ws.Range("I134").Value = 2366.3076
$ws.Range("J134").Value = 1241.909
$ws.Range("K134").Value = 7098.9228
$ws.Range("L134").Value = 3725.727
$ws.Range("M134").Value = -4563.9228
$ws.Range("N134").Value = -8795.727

$ws = $wb.Worksheets.Item("CRP")
# Row 36
$ws.Range("H36").Value = 0
$ws.Range("I36").Value = 0
$ws.Range("K36").Value = 0
$ws.Range("M36").ClearContents()
# Row 40
$ws.Range("H40").Value = 0
$ws.Range("I40").Value = 0
$ws.Range("K40").Value = 0
$ws.Range("M40").ClearContents()
# Row 86
$ws.Range("H86").Value = 3730.7646
$ws.Range("I86").Value = 3580.1428
$ws.Range("K86").Value = 3580.1428
$ws.Range("M86").Value = -2457.1428
# Row 89
$ws.Range("H89").Value = 3730.7646
$ws.Range("I89").Value = 3580.1428
$ws.Range("K89").Value = 17900.714
$ws.Range("M89").Value = -12284.714
# Row 122
$ws.Range("H122").Value = 1392.6897
$ws.Range("I122").Value = 1042.2222
$ws.Range("J122").Value = 1966.1818
$ws.Range("K122").Value = 3126.6666
$ws.Range("L122").Value = 5898.5454
$ws.Range("M122").Value = -676.6665999999996
$ws.Range("N122").Value = -10798.5454

$ws = $wb.Worksheets.Item("CUL")
# Row 3
$ws.Range("H3").Value = 2821.1765
$ws.Range("I3").Value = 2821.1765
$ws.Range("K3").Value = 8463.5295
$ws.Range("M3").Value = -8351.5295
# Row 43
$ws.Range("H43").Value = 4140.4
$ws.Range("I43").Value = 2802
$ws.Range("K43").Value = 8406
$ws.Range("M43").Value = -8292
# Row 113
$ws.Range("H113").Value = 533.2308
$ws.Range("I113").Value = 493.04544
$ws.Range("J113").Value = 585.2353
$ws.Range("K113").Value = 1479.13632
$ws.Range("L113").Value = 1755.7059
$ws.Range("M113").Value = 690.8636799999999
$ws.Range("N113").Value = -6095.7059
# Row 133
$ws.Range("H133").Value = 3404.1904
$ws.Range("I133").Value = 2638.3333
$ws.Range("J133").Value = 7999.3335
$ws.Range("K133").Value = 7914.999899999999
$ws.Range("L133").Value = 23998.0005
$ws.Range("M133").Value = -2854.999899999999
$ws.Range("N133").Value = -34118.00049999999

$ws = $wb.Worksheets.Item("GSM")
# Row 80
$ws.Range("H80").Value = 7388.08
$ws.Range("I80").Value = 11554.546
$ws.Range("J80").Value = 4114.4287
$ws.Range("K80").Value = 11554.546
$ws.Range("L80").Value = 4114.4287
$ws.Range("M80").Value = -10556.546
$ws.Range("N80").Value = -6110.4287
# Row 83
$ws.Range("H83").Value = 7388.08
$ws.Range("I83").Value = 11554.546
$ws.Range("J83").Value = 4114.4287
$ws.Range("K83").Value = 57772.73
$ws.Range("L83").Value = 20572.1435
$ws.Range("M83").Value = -52780.73
$ws.Range("N83").Value = -30556.1435
# Row 125
$ws.Range("H125").Value = 44596.668
$ws.Range("J125").Value = 44596.668
$ws.Range("L125").Value = 44596.668
$ws.Range("N125").Value = -49516.668

$ws = $wb.Worksheets.Item("LTW")
# Row 39
$ws.Range("H39").Value = 25000
$ws.Range("J39").Value = 25000
$ws.Range("L39").Value = 25000
$ws.Range("N39").Value = -25920
# Row 93
$ws.Range("H93").Value = 1271.4286
$ws.Range("I93").Value = 1340
$ws.Range("J93").Value = 1100
$ws.Range("K93").Value = 1340
$ws.Range("L93").Value = 1100
$ws.Range("M93").Value = -92
$ws.Range("N93").Value = -3596
# Row 100
$ws.Range("H100").Value = 34733.332
$ws.Range("I100").Value = 67766.664
$ws.Range("K100").Value = 67766.664
$ws.Range("M100").Value = -67225.664

$ws = $wb.Worksheets.Item("WVR")
# Row 135
$ws.Range("H135").Value = 43428.5
$ws.Range("J135").Value = 43428.5
$ws.Range("L135").Value = 43428.5
$ws.Range("N135").Value = -53568.5
